# Generate Report for Handback
#
# - "Status" for the 1f59c2f2-... row flips from "Ready for handoff" to
#   "Handback transform failed" on every sheet that shows it (Overview's
#   zh-cn/de-de columns and the per-locale Status column).
# - The per-locale "Error Detail" cell (column P) for that same row gets a
#   new diagnostic message explaining the handback/handoff filename
#   mismatch.
# - The "Error Detail" column is widened to fit the new message.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = "Handback file name: 4dyxnl1m.ir1 is different with handoff file name: 1f59c2f2-e1dd-463f-ab74-8bdc62288281.22d0a7fe9c3f3190fe85f57c41263c570b36e9b2.zh-cn."
$zhcn.Range("P1").ColumnWidth = 39.2

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = "Handback file name: 4dyxnl1m.ir1 is different with handoff file name: 1f59c2f2-e1dd-463f-ab74-8bdc62288281.22d0a7fe9c3f3190fe85f57c41263c570b36e9b2.de-de."
$dede.Range("P1").ColumnWidth = 39.2
